$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 466850.4
$ws.Range("J17").Value = 466850.4
$ws.Range("L17").Value = 1400551.2
$ws.Range("N17").Value = -1400887.2
$ws.Range("H18").Value = 19813.834
$ws.Range("I18").Value = 22960.2
$ws.Range("K18").Value = 22960.2
$ws.Range("M18").Value = -22676.2
$ws.Range("H113").Value = 147857.14
$ws.Range("I113").Value = 500500
$ws.Range("J113").Value = 6800
$ws.Range("K113").Value = 500500
$ws.Range("L113").Value = 6800
$ws.Range("M113").Value = -497246
$ws.Range("N113").Value = -13308
$ws.Range("H138").Value = 2229.2327
$ws.Range("I138").Value = 1313.4762
$ws.Range("J138").Value = 3103.3635
$ws.Range("K138").Value = 3940.4286
$ws.Range("L138").Value = 9310.0905
$ws.Range("M138").Value = 1199.5714
$ws.Range("N138").Value = -19590.0905

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1624.4166
$ws.Range("I2").Value = 1499.5428
$ws.Range("K2").Value = 1499.5428
$ws.Range("M2").Value = -1386.5428
$ws.Range("H37").Value = 24994.8
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 24994.8
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 24994.8
$ws.Range("N37").Value = -25540.8
$ws.Range("H45").Value = 9218.058999999999
$ws.Range("I45").Value = 18584.5
$ws.Range("K45").Value = 18584.5
$ws.Range("M45").Value = -18207.5
$ws.Range("H61").Value = 1778.973
$ws.Range("I61").Value = 1689.5
$ws.Range("K61").Value = 1689.5
$ws.Range("M61").Value = -1477.5
$ws.Range("H74").Value = 1284.9286
$ws.Range("I74").Value = 1133
$ws.Range("K74").Value = 1133
$ws.Range("M74").Value = -259
$ws.Range("H77").Value = 1284.9286
$ws.Range("I77").Value = 1133
$ws.Range("K77").Value = 5665
$ws.Range("M77").Value = -1297
$ws.Range("H116").Value = 1624.4166
$ws.Range("I116").Value = 1499.5428
$ws.Range("K116").Value = 1499.5428
$ws.Range("M116").Value = 794.4572000000001
$ws.Range("H122").Value = 1885.0889
$ws.Range("I122").Value = 1510.7931
$ws.Range("K122").Value = 4532.379300000001
$ws.Range("M122").Value = -2082.379300000001
$ws.Range("H132").Value = 5020.9165
$ws.Range("I132").Value = 3531.375
$ws.Range("K132").Value = 10594.125
$ws.Range("M132").Value = -8064.125
$ws.Range("H136").Value = 1778.973
$ws.Range("I136").Value = 1689.5
$ws.Range("K136").Value = 5068.5
$ws.Range("M136").Value = -2518.5
$ws.Range("M37").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1624.4166
$ws.Range("I3").Value = 1499.5428
$ws.Range("K3").Value = 1499.5428
$ws.Range("M3").Value = -1385.5428
$ws.Range("H94").Value = 1095.1538
$ws.Range("J94").Value = 744.5
$ws.Range("L94").Value = 744.5
$ws.Range("N94").Value = -1646.5
$ws.Range("H105").Value = 1658.7858
$ws.Range("I105").Value = 1563.3077
$ws.Range("K105").Value = 1563.3077
$ws.Range("M105").Value = 183.6922999999999
$ws.Range("H134").Value = 3180.6667
$ws.Range("I134").Value = 2861.3333
$ws.Range("K134").Value = 8583.999899999999
$ws.Range("M134").Value = -6048.999899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1140.3
$ws.Range("I22").Value = 600.3333
$ws.Range("J22").Value = 1371.7142
$ws.Range("K22").Value = 600.3333
$ws.Range("L22").Value = 1371.7142
$ws.Range("M22").Value = -250.3333
$ws.Range("N22").Value = -2071.7142
$ws.Range("H31").Value = 3017.6943
$ws.Range("I31").Value = 1904.1333
$ws.Range("K31").Value = 1904.1333
$ws.Range("H34").Value = 3017.6943
$ws.Range("I34").Value = 1904.1333
$ws.Range("K34").Value = 1904.1333
$ws.Range("H44").Value = 4500
$ws.Range("I44").Value = 4500
$ws.Range("K44").Value = 4500
$ws.Range("H55").Value = 10000
$ws.Range("J55").Value = 10000
$ws.Range("L55").Value = 10000
$ws.Range("H132").Value = 4338.636
$ws.Range("J132").Value = 2385.4285
$ws.Range("L132").Value = 7156.2855
$ws.Range("N132").Value = -12216.2855
$ws.Range("H134").Value = 1911.9736
$ws.Range("I134").Value = 1078.1177
$ws.Range("K134").Value = 3234.3531
$ws.Range("M134").Value = -699.3531000000003
$ws.Range("H141").Value = 376330.34
$ws.Range("J141").Value = 376330.34
$ws.Range("L141").Value = 376330.34
$ws.Range("N141").Value = -386690.34
$ws.Range("M31").Value = -1609.1333
$ws.Range("M34").Value = -1702.1333
$ws.Range("M44").Value = -4058
$ws.Range("N55").Value = -10630

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 9825.532999999999
$ws.Range("I56").Value = 9825.532999999999
$ws.Range("K56").Value = 9825.532999999999
$ws.Range("M56").Value = -9295.532999999999
$ws.Range("H86").Value = 993.375
$ws.Range("I86").Value = 974
$ws.Range("K86").Value = 2922
$ws.Range("M86").Value = -1736
$ws.Range("H89").Value = 993.375
$ws.Range("I89").Value = 974
$ws.Range("K89").Value = 8766
$ws.Range("M89").Value = -2838
$ws.Range("H107").Value = 593.41174
$ws.Range("I107").Value = 826.6667
$ws.Range("J107").Value = 570.8387
$ws.Range("K107").Value = 2480.0001
$ws.Range("L107").Value = 1712.5161
$ws.Range("M107").Value = -560.0001000000002
$ws.Range("N107").Value = -5552.5161
$ws.Range("H128").Value = 468476
$ws.Range("I128").Value = 468476
$ws.Range("K128").Value = 1405428
$ws.Range("M128").Value = -1400448
$ws.Range("H131").Value = 2366.2593
$ws.Range("J131").Value = 2727.1177
$ws.Range("L131").Value = 8181.353099999999
$ws.Range("N131").Value = -18261.3531
$ws.Range("H137").Value = 3032.8918
$ws.Range("I137").Value = 1633.3636
$ws.Range("J137").Value = 3625
$ws.Range("K137").Value = 4900.0908
$ws.Range("L137").Value = 10875
$ws.Range("M137").Value = 199.9092000000001
$ws.Range("N137").Value = -21075

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3104
$ws.Range("I70").Value = 3104
$ws.Range("K70").Value = 3104
$ws.Range("H73").Value = 3104
$ws.Range("I73").Value = 3104
$ws.Range("K73").Value = 3104
$ws.Range("H102").Value = 1339.7715
$ws.Range("J102").Value = 831.3333
$ws.Range("L102").Value = 831.3333
$ws.Range("N102").Value = -4075.3333
$ws.Range("H122").Value = 2459.7307
$ws.Range("J122").Value = 1785.2
$ws.Range("L122").Value = 5355.6
$ws.Range("N122").Value = -10255.6
$ws.Range("M70").Value = -2834
$ws.Range("M73").Value = -2168

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I93").Value = 6415
$ws.Range("J93").Value = 52440.855
$ws.Range("K93").Value = 6415
$ws.Range("L93").Value = 52440.855
$ws.Range("M93").Value = -5167
$ws.Range("N93").Value = -54936.855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1628.4615
$ws.Range("I122").Value = 1495.8948
$ws.Range("K122").Value = 4487.6844
$ws.Range("M122").Value = -2037.6844
$ws.Range("H126").Value = 2775.5
$ws.Range("I126").Value = 2730.8
$ws.Range("K126").Value = 8192.400000000001
$ws.Range("M126").Value = -5722.400000000001
